$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 54,4
$data[0,0] = -0.024166140465835009
$data[0,1] = 0.028113961550550574
$data[0,2] = -0.079269173814557309
$data[0,3] = 0.030936892882887288
$data[1,0] = -0.055181831103036753
$data[1,1] = 0.018648370263095789
$data[1,2] = -0.091732302474690608
$data[1,3] = -0.018631359731382899
$data[2,0] = -0.036284606564157322
$data[2,1] = 0.033797935233883913
$data[2,2] = -0.10252840091460541
$data[2,3] = 0.029959187786290774
$data[3,0] = -0.058452591592142375
$data[3,1] = 0.022357307828168478
$data[3,2] = -0.10227260914933772
$data[3,3] = -0.01463257403494702
$data[4,0] = -0.16500529013370988
$data[4,1] = 0.077922926397031175
$data[4,2] = -0.31773989285240317
$data[4,3] = -0.012270687415016579
$data[5,0] = -0.091196460510440155
$data[5,1] = 0.083541725275067286
$data[5,2] = -0.25494329907179519
$data[5,3] = 0.072550378050914868
$data[6,0] = -0.091935351397948062
$data[6,1] = 0.030762047325026898
$data[6,2] = -0.15222860166003127
$data[6,3] = -0.031642101135864853
$data[7,0] = -0.067358595908313154
$data[7,1] = 0.023097571922626389
$data[7,2] = -0.1126294227635347
$data[7,3] = -0.022087769053091608
$data[8,0] = -0.068268189828735282
$data[8,1] = 0.030915524094013092
$data[8,2] = -0.12886247188023039
$data[8,3] = -0.0076739077772401687
$data[9,0] = -0.047445227704639457
$data[9,1] = 0.020880982780650158
$data[9,2] = -0.088371668360734387
$data[9,3] = -0.0065187870485445276
$data[10,0] = -0.048241053036511682
$data[10,1] = 0.069550140967599233
$data[10,2] = -0.18456438739886316
$data[10,3] = 0.088082281325839784
$data[11,0] = -0.065559608273461509
$data[11,1] = 0.065128809607767277
$data[11,2] = -0.19321601753987533
$data[11,3] = 0.062096800992952297
$data[12,0] = -0.090418762784684234
$data[12,1] = 0.026430590720041034
$data[12,2] = -0.14222240914217629
$data[12,3] = -0.038615116427192168
$data[13,0] = -0.085830021606105533
$data[13,1] = 0.021438927481730565
$data[13,2] = -0.12784993509477471
$data[13,3] = -0.043810108117436346
$data[14,0] = -0.076606682476943286
$data[14,1] = 0.035197824354045171
$data[14,2] = -0.14559425292932385
$data[14,3] = -0.0076191120245627242
$data[15,0] = -0.091519529243567752
$data[15,1] = 0.021988850096528078
$data[15,2] = -0.13461737468622864
$data[15,3] = -0.048421683800906876
$data[16,0] = -0.054582204079242289
$data[16,1] = 0.065214203243786031
$data[16,2] = -0.18240678516986297
$data[16,3] = 0.07324237701137841
$data[17,0] = -0.074148744668780256
$data[17,1] = 0.052170263493979602
$data[17,2] = -0.17640561913389946
$data[17,3] = 0.028108129796338949
$data[18,0] = -0.056494154079221641
$data[18,1] = 0.025345225757008861
$data[18,2] = -0.1061704978989351
$data[18,3] = -0.0068178102595081932
$data[19,0] = -0.063963121872785494
$data[19,1] = 0.022855721626186692
$data[19,2] = -0.10875992648308527
$data[19,3] = -0.019166317262485723
$data[20,0] = -0.076902034852655385
$data[20,1] = 0.030989975147923617
$data[20,2] = -0.13764224062020872
$data[20,3] = -0.01616182908510206
$data[21,0] = -0.09003152329781508
$data[21,1] = 0.025037535155212381
$data[21,2] = -0.13910474975793405
$data[21,3] = -0.040958296837696105
$data[22,0] = -0.084479695707142116
$data[22,1] = 0.057092958198765024
$data[22,2] = -0.19638604588825609
$data[22,3] = 0.027426654473971843
$data[23,0] = -0.0074946105297094131
$data[23,1] = 0.055440832646043381
$data[23,2] = -0.11616199850988247
$data[23,3] = 0.10117277745046364
$data[24,0] = -0.033393313471913932
$data[24,1] = 0.030261433495671501
$data[24,2] = -0.092705366527612632
$data[24,3] = 0.02591873958378476
$data[25,0] = -0.080148514271386029
$data[25,1] = 0.021310138515043882
$data[25,2] = -0.12191600369438862
$data[25,3] = -0.038381024848383426
$data[26,0] = -0.066241074821756563
$data[26,1] = 0.031636135876366982
$data[26,2] = -0.12824775258282459
$data[26,3] = -0.0042343970606885314
$data[27,0] = -0.051391113130641021
$data[27,1] = 0.024991057587349153
$data[27,2] = -0.10037324419343241
$data[27,3] = -0.0024089820678496396
$data[28,0] = -0.031667060526725931
$data[28,1] = 0.042201270167320976
$data[28,2] = -0.11438461916137682
$data[28,3] = 0.051050498107924967
$data[29,0] = 0.0089077851892962294
$data[29,1] = 0.048914563461429736
$data[29,2] = -0.08696772013508719
$data[29,3] = 0.10478329051367966
$data[30,0] = -0.01109517491466747
$data[30,1] = 0.030845134577347902
$data[30,2] = -0.071551275212213161
$data[30,3] = 0.049360925382878215
$data[31,0] = -0.06573425887274513
$data[31,1] = 0.022894782335718646
$data[31,2] = -0.11060762177341354
$data[31,3] = -0.020860895972076708
$data[32,0] = -0.026938170760072652
$data[32,1] = 0.030472520558056118
$data[32,2] = -0.086664167961343117
$data[32,3] = 0.03278782644119782
$data[33,0] = -0.091459032832886478
$data[33,1] = 0.027344536422461006
$data[33,2] = -0.14505395022342854
$data[33,3] = -0.037864115442344401
$data[34,0] = -0.0444606225452591
$data[34,1] = 0.046305461116525448
$data[34,2] = -0.13522269391951733
$data[34,3] = 0.046301448828999135
$data[35,0] = -0.060454963080329459
$data[35,1] = 0.041688415827932181
$data[35,2] = -0.14216678162344429
$data[35,3] = 0.021256855462785378
$data[36,0] = -0.0031838329449517546
$data[36,1] = 0.028484624903416085
$data[36,2] = -0.059013362097453705
$data[36,3] = 0.052645696207550198
$data[37,0] = -0.039906331028275781
$data[37,1] = 0.026155540925570247
$data[37,2] = -0.091170722303370103
$data[37,3] = 0.011358060246818541
$data[38,0] = -0.042458771893095472
$data[38,1] = 0.036840130664480178
$data[38,2] = -0.11466525500200514
$data[38,3] = 0.029747711215814189
$data[39,0] = -0.074145767893833966
$data[39,1] = 0.031778423784631145
$data[39,2] = -0.13643104387093752
$data[39,3] = -0.011860491916730415
$data[40,0] = -0.1222383702061544
$data[40,1] = 0.038990697280843763
$data[40,2] = -0.1986629724929263
$data[40,3] = -0.045813767919382492
$data[41,0] = -0.013034892039827019
$data[41,1] = 0.044961169489526082
$data[41,2] = -0.10116150586890943
$data[41,3] = 0.075091721789255383
$data[42,0] = 0.000076107635260761524
$data[42,1] = 0.035965042230396767
$data[42,2] = -0.07041495133010045
$data[42,3] = 0.070567166600621983
$data[43,0] = -0.065241118925153455
$data[43,1] = 0.027624462402308331
$data[43,2] = -0.11938456995856736
$data[43,3] = -0.011097667891739542
$data[44,0] = -0.0096562042798023842
$data[44,1] = 0.038657440026109401
$data[44,2] = -0.08542460520380446
$data[44,3] = 0.066112196644199678
$data[45,0] = -0.074803561153062359
$data[45,1] = 0.030208396366586027
$data[45,2] = -0.13401160486442171
$data[45,3] = -0.015595517441702997
$data[46,0] = -0.2176558493150007
$data[46,1] = 0.049961419375243198
$data[46,2] = -0.3155838647575846
$data[46,3] = -0.11972783387241683
$data[47,0] = -0.026113855041104767
$data[47,1] = 0.054385118238539254
$data[47,2] = -0.13271197887712366
$data[47,3] = 0.080484268794914132
$data[48,0] = -0.060820636204932041
$data[48,1] = 0.041544280060603012
$data[48,2] = -0.14224693557266505
$data[48,3] = 0.020605663162800977
$data[49,0] = -0.080862364669894687
$data[49,1] = 0.037395867423835855
$data[49,2] = -0.15415759435535292
$data[49,3] = -0.0075671349844364383
$data[50,0] = 0.0037547668028495526
$data[50,1] = 0.058984707574665556
$data[50,2] = -0.11185498306376809
$data[50,3] = 0.11936451666946721
$data[51,0] = -0.12040337143116844
$data[51,1] = 0.041324012088176473
$data[51,2] = -0.20139786992602543
$data[51,3] = -0.039408872936311456
$data[52,0] = -0.09621075419213021
$data[52,1] = 0.06933949645387133
$data[52,2] = -0.23212120998835625
$data[52,3] = 0.039699701604095833
$data[53,0] = 0.040782263691364777
$data[53,1] = 0.073856337348246412
$data[53,2] = -0.10398062824888119
$data[53,3] = 0.18554515563161075
$ws.Range("B2:E55").Value = $data
